$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:F13"))

$ws.Range("F12").Value = "TC212B106K020Y Sunlord | C177761 - LCSC Electronics"
$ws.Range("B12").Value = "Tantalum capacitor"
$ws.Range("F13").Value = "TAJB226K010RNJ Kyocera AVX | C7198 - LCSC Electronics"
$ws.Range("A13").Value = "TAJB226K010RNJ"
$ws.Range("A12").Value = "TC212B106K020Y"

$ws.Range("C12").Value = 0.24
$ws.Range("D12").Value = 1
$ws.Range("E12").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

$ws.Range("B13").Value = "Tantalum capacitor"
$ws.Range("C13").Value = 0.07
$ws.Range("D13").Value = 1
$ws.Range("E13").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lcsc.com/product-detail/Tantalum-Capacitors_Sunlord-TC212B106K020Y_C177761.html", "", "", "https://www.lcsc.com/product-detail/Tantalum-Capacitors_Sunlord-TC212B106K020Y_C177761.html")
$ws.Range("F12").Value = "TC212B106K020Y Sunlord | C177761 - LCSC Electronics"
$ws.Range("F12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lcsc.com/product-detail/Tantalum-Capacitors_Kyocera-AVX-TAJB226K010RNJ_C7198.html", "", "", "https://www.lcsc.com/product-detail/Tantalum-Capacitors_Kyocera-AVX-TAJB226K010RNJ_C7198.html")
$ws.Range("F13").Value = "TAJB226K010RNJ Kyocera AVX | C7198 - LCSC Electronics"
$ws.Range("F13").Style = "Hyperlink"

[void]$ws.Range("B17").Select()
